$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.4308180809021
$ws.Range("B1").Value = 4.401615142822266
$ws.Range("C1").Value = 7.271471500396729
$ws.Range("D1").Value = 8.037375450134277
$ws.Range("E1").Value = 5.783596515655518
